$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.170.94"
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = "  -0.15%  "

$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.903.19"
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = "  +0.11%  "

$ws.Range("E4").Value = "  -0.18%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.18"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +2.84%  "

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.695"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("E7").Value = "  -0.06%  "

$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.30"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  -1.15%  "

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.362"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +3.45%  "

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.96"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -1.20%  "

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0754"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +3.63%  "

$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0983"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  -1.15%  "

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.14"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  +6.65%  "

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.182.57"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +0.12%  "

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.740"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  +4.80%  "

$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.00"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  +3.14%  "

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.906.90"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  +0.21%  "

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.191.96"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  -0.30%  "

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.03"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +2.24%  "

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0836"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +1.53%  "

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "242.52"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  +0.55%  "

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.97"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +3.31%  "

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.07"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +4.76%  "

$ws.Range("E24").Value = "  -0.15%  "

$ws.Range("E25").Value = "  +5.62%  "

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.31"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  -1.26%  "

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.31"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -1.99%  "

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.63"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +0.89%  "

$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.50"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +0.54%  "

$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.130"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  -1.05%  "

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.127.88"
$ws.Range("D31").Style = $origStyle

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.13"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  +19.07%  "

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0621"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  +8.85%  "

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.32"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +3.10%  "

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.59"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  +18.01%  "

$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.19"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  +1.70%  "

$ws.Range("E37").Value = "  -0.11%  "

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.858"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  -12.62%  "

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.01"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  -1.68%  "

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "100.01"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +10.18%  "

$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.26"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +5.94%  "

$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0215"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +2.37%  "

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.11"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  +0.86%  "

$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0652"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -3.76%  "

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.46"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  +1.74%  "

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.335.80"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -0.60%  "

$ws.Range("E47").Value = "  +0.59%  "

$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.75"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -1.41%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.64"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  +1.13%  "

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.87"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -7.98%  "

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.18"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -8.98%  "
